$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting the existing weekly records (rows 6-19)
# down to rows 7-20. This mirrors the prior week's data being appended as the
# most recent entry while all older rows move down one position.
$ws.Rows.Item(6).Insert()

# Preserve the date-formatted style used by the other rows in column D.
$ws.Range("D6").NumberFormat = $ws.Range("D7").NumberFormat

# Populate the new row 6 with this week's record.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44838
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112012
$ws.Range("G6").Value = "Espinaca"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6750
$ws.Range("N6").Value = "$/cuna 10 kilos"
$ws.Range("O6").Value = "Provincia de Diguillín"
$ws.Range("P6").Value = 675
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
